# Update countries & provincias Spain
# Refreshes the "Pais" sheet with the latest COVID-19 case counts and
# re-sorts a handful of country rows whose totals changed rank order,
# plus updates the "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 23 de Junio de 2020 a las 00:35"

# Row 4
$ws.Cells.Item(4, 2).Value = 2385905
$ws.Cells.Item(4, 3).Value = 29248
$ws.Cells.Item(4, 4).Value = 990261
$ws.Cells.Item(4, 5).Value = 1273055
$ws.Cells.Item(4, 7).Value = 342
$ws.Cells.Item(4, 8).Value = 122589

# Row 5
$ws.Cells.Item(5, 2).Value = 1106470
$ws.Cells.Item(5, 3).Value = 19480
$ws.Cells.Item(5, 5).Value = 475973
$ws.Cells.Item(5, 7).Value = 612
$ws.Cells.Item(5, 8).Value = 51271

# Row 14
$ws.Cells.Item(14, 2).Value = 192074
$ws.Cells.Item(14, 3).Value = 499
$ws.Cells.Item(14, 5).Value = 7805
$ws.Cells.Item(14, 7).Value = 7
$ws.Cells.Item(14, 8).Value = 8969

# Row 19
$ws.Cells.Item(19, 4).Value = 74612
$ws.Cells.Item(19, 5).Value = 56475

# Row 25
$ws.Cells.Item(25, 2).Value = 71183
$ws.Cells.Item(25, 3).Value = 2531
$ws.Cells.Item(25, 4).Value = 28968
$ws.Cells.Item(25, 5).Value = 39905
$ws.Cells.Item(25, 7).Value = 73
$ws.Cells.Item(25, 8).Value = 2310

# Row 28
$ws.Cells.Item(28, 1).Value = "Suecia"
$ws.Cells.Item(28, 2).Value = 58932
$ws.Cells.Item(28, 3).Value = 84
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 7).Value = 69
$ws.Cells.Item(28, 8).Value = 5122

# Row 29
$ws.Cells.Item(29, 1).Value = "Egipto"
$ws.Cells.Item(29, 2).Value = 56809
$ws.Cells.Item(29, 3).Value = 1576
$ws.Cells.Item(29, 4).Value = 15133
$ws.Cells.Item(29, 5).Value = 39398
$ws.Cells.Item(29, 7).Value = 85
$ws.Cells.Item(29, 8).Value = 2278

# Row 52
$ws.Cells.Item(52, 1).Value = "Nigeria"
$ws.Cells.Item(52, 2).Value = 20919
$ws.Cells.Item(52, 3).Value = 675
$ws.Cells.Item(52, 4).Value = 7109
$ws.Cells.Item(52, 5).Value = 13285
$ws.Cells.Item(52, 7).Value = 7
$ws.Cells.Item(52, 8).Value = 525

# Row 53
$ws.Cells.Item(53, 1).Value = "Armenia"
$ws.Cells.Item(53, 2).Value = 20588
$ws.Cells.Item(53, 3).Value = 320
$ws.Cells.Item(53, 4).Value = 9131
$ws.Cells.Item(53, 5).Value = 11097
$ws.Cells.Item(53, 7).Value = 10
$ws.Cells.Item(53, 8).Value = 360

# Row 90
$ws.Cells.Item(90, 1).Value = "Bulgaria"
$ws.Cells.Item(90, 2).Value = 3984
$ws.Cells.Item(90, 3).Value = 79
$ws.Cells.Item(90, 4).Value = 2171
$ws.Cells.Item(90, 5).Value = 1606
$ws.Cells.Item(90, 7).Value = 8
$ws.Cells.Item(90, 8).Value = 207

# Row 91
$ws.Cells.Item(91, 1).Value = "Venezuela"
$ws.Cells.Item(91, 2).Value = 3917
$ws.Cells.Item(91, 4).Value = 835
$ws.Cells.Item(91, 5).Value = 3049
$ws.Cells.Item(91, 8).Value = 33

# Row 167
$ws.Cells.Item(167, 4).Value = 146
$ws.Cells.Item(167, 5).Value = 48

# Row 190
$ws.Cells.Item(190, 1).Value = "Gambia"
$ws.Cells.Item(190, 2).Value = 41
$ws.Cells.Item(190, 3).Value = 4
$ws.Cells.Item(190, 4).Value = 26
$ws.Cells.Item(190, 5).Value = 13

# Row 191
$ws.Cells.Item(191, 1).Value = "Puerto Rico"
$ws.Cells.Item(191, 2).Value = 39
$ws.Cells.Item(191, 4).Value = 1
$ws.Cells.Item(191, 5).Value = 36

# Row 202
$ws.Cells.Item(202, 1).Value = "Fiyi"

# Row 203
$ws.Cells.Item(203, 1).Value = "Dominica"
